# The site's footer block ("Ver no Jupiter Salvar em pdf Salvar em docx",
# the copyright/contact line, and the blank paragraph right after it) was
# dropped from the page during a site rebuild. Remove that block, leaving
# the surrounding (already-blank) paragraphs — including the trailing
# page-break paragraph — untouched.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count

$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text

    if ($paraText -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startIndex = $i
    }

    if ($paraText -like "*Contact: luizeleno@usp.br*") {
        # also swallow the blank paragraph immediately following the
        # copyright line
        $endIndex = $i + 1
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $startRange = $d.Paragraphs.Item($startIndex).Range
    $endRange = $d.Paragraphs.Item($endIndex).Range

    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}

Write-Output $d.Paragraphs.Count
